# Applies the CCDI TC01 input-file touch-ups captured in the commit
# "ccdi - small updates to the input files based on xl macro":
#   - SamplesTab / StudiesTab query cells are rewritten (same text, but the
#     rewrite order changes how the shared-string table is built on save).
#   - DiagnosisTab query text: "Age at diagnosis (days)" -> "Age at Diagnosis (days)".
#   - FilesTab query text: "MD5Sum" -> "MD5sum".
#   - Active selection on the startup sheet moves from B4 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rewrite order matters: it reproduces the shared-string allocation order
# (Samples, Studies, Diagnosis, Files) seen in the target workbook.
$ws.Range('B5').Value2 = 'CALL {
  MATCH (st:study)<-[:of_participant|of_cell_line|of_pdx]-(pcp)<-[:of_sample]-(sm:sample)
  WHERE (pcp:participant or pcp:cell_line or pcp:pdx) and st.phs_accession IN [''phs003111'']
  return st, case labels(pcp)[0] when "participant" then pcp.participant_id else "" end as participant_id, sm
  union all
  MATCH (st:study)<-[:of_participant]-(p:participant)<-[:of_sample]-(sm1:sample)<--(cp)<--(sm2:sample)
  WHERE (cp:cell_line or cp:pdx) and st.phs_accession IN [''phs003111'']
  with st, p.participant_id as participant_id, [sm1,sm2] as sampleArr
  unwind sampleArr as sm
  return st, participant_id, sm
}
RETURN DISTINCT
          sm.sample_id as `Sample ID`,
          participant_id as `Participant ID`,
          st.study_id as `Study ID`,
          sm.anatomic_site as `Anatomic Site`,
          case sm.participant_age_at_collection when -999 then ''Not Reported'' else coalesce(sm.participant_age_at_collection, '''') end as `Age at Sample Collection`,
          coalesce(sm.diagnosis_classification, '''') as `Diagnosis`,
          coalesce(sm.diagnosis_classification_system, '''') as `Diagnosis Classification System`,
          coalesce(sm.diagnosis_verification_status, '''') as `Diagnosis Verification Status`,
          coalesce(sm.diagnosis_basis, '''') as `Diagnosis Basis`,
          coalesce(sm.diagnosis_comment, '''') as `Diagnosis Comment`,
          sm.sample_tumor_status as `Sample Tumor Status`,
          sm.tumor_classification as `Sample Tumor Classification`
Order by sm.sample_id Limit 100'

$ws.Range('B4').Value2 = 'MATCH (st:study)<-[:of_participant]-(p:participant)
where st.phs_accession in ["phs003111"]
with st, count(p) as num_p
MATCH (st)<-[:of_participant]-(participant)<-[:of_diagnosis]-(dg:diagnosis)
with st, num_p, dg.diagnosis_classification as dg_cancers, count(dg.diagnosis_classification) as num_cancers
ORDER BY num_cancers desc
with st, num_p, collect(dg_cancers + '' ('' + toString(num_cancers) + '')'') as cancers
MATCH (st)<-[:of_participant]-(pa:participant)<-[:of_diagnosis]-(diag:diagnosis)
with st, num_p, cancers, diag.anatomic_site as dg_sites, count(diag.anatomic_site) as num_sites
ORDER BY num_sites desc
with st, num_p, cancers, collect(dg_sites + '' ('' + toString(num_sites) + '')'') as sites
MATCH (st)<-[*..5]-(fl)
WHERE (fl:clinical_measure_file OR fl: sequencing_file OR fl:pathology_file OR fl:radiology_file OR fl:methylation_array_file OR fl:single_cell_sequencing_file OR fl:cytogenomic_file)
with st, num_p, cancers, sites, fl.file_type as ft, count(fl.file_type) as num_ft
ORDER BY num_ft desc
with st, num_p, cancers, sites, collect(ft + '' ('' + toString(num_ft) + '')'') as file_types, sum(num_ft) as num_files
OPTIONAL MATCH (st)<-[:of_participant|of_cell_line|of_pdx]-(pcp)<-[:of_sample]-(sm1:sample)
WHERE (pcp:participant or pcp:cell_line or pcp:pdx)
WITH st, num_p, cancers, sites, file_types, num_files, count(distinct sm1.sample_id) as num_samples_1
OPTIONAL MATCH (st)<-[:of_participant]-(participant)<-[:of_sample]-(sm1:sample)<--(cp)<--(sm2:sample)
WHERE (cp:cell_line or cp:pdx)
WITH st, num_p, cancers, sites, file_types, num_files, num_samples_1, count(distinct sm2.sample_id) as num_samples_2
WITH st, num_p, cancers, sites, file_types, num_files, num_samples_1 + num_samples_2 as num_samples
MATCH (st)<-[*..5]-(file)
WHERE (file:clinical_measure_file OR file: sequencing_file OR file:pathology_file OR file:radiology_file OR file:methylation_array_file OR file:single_cell_sequencing_file OR file:cytogenomic_file)
OPTIONAL MATCH (st)<-[:of_publication]-(pub:publication)
OPTIONAL MATCH (st)<-[:of_study_personnel]-(stp:study_personnel)
WHERE stp.personnel_type = ''PI''
OPTIONAL MATCH (st)<-[:of_study_funding]-(stf:study_funding)
WITH st, num_p, cancers, sites, file_types, num_files, num_samples, file.id as file_id, stf, stp, pub
order by st.study_id
RETURN DISTINCT
  st.study_short_title as `Study Short Title`,
st.study_id as `Study ID`,
  CASE WHEN size(cancers) > 5 THEN apoc.text.join(apoc.coll.remove(cancers, 5, 10000), "") + "Read More"  else apoc.text.join(cancers, "") END as `Diagnosis (Top 5)`,
  CASE WHEN size(sites) > 5 THEN apoc.text.join(apoc.coll.remove(sites, 5, 10000), "") + "Read More"  else apoc.text.join(sites, "") END as `Diagnosis Anatomic Site (Top 5)`,
  num_p as `Number of Participants`,
  num_samples as `Number of Samples`,
  num_files as `Number of Files`,
  CASE WHEN size(file_types) > 5 THEN apoc.text.join(apoc.coll.remove(file_types, 5, 10000), "") + "Read More"  else apoc.text.join(file_types, "") END as `File Type (Top 5)`,
  apoc.text.join(COLLECT(DISTINCT pub.pubmed_id), '';'') as `PubMed ID`,
  apoc.text.join(COLLECT(DISTINCT stp.personnel_name), '';'') as `Principal Investigator(s)`,
  apoc.text.join(COLLECT(DISTINCT stf.grant_id), '';'') as `Grant ID`'

$ws.Range('B3').Value2 = 'Match (s:study)<--(p:participant)<--(d:diagnosis)
where s.phs_accession = "phs003111"
OPTIONAL MATCH(fo:follow_up)-->(p)
with distinct d, p, s, fo
return
coalesce(p.participant_id, '''') as `Participant ID`,
coalesce(s.phs_accession, '''') as `Study ID`,
coalesce(d.diagnosis_classification, '''') as `Diagnosis`,
coalesce(d.diagnosis_classification_system, '''') as `Diagnosis Classification System`,
coalesce(d.diagnosis_verification_status, '''') as `Diagnosis Verification Status`,
coalesce(d.diagnosis_basis, '''') as `Diagnosis Basis`,
coalesce(d.diagnosis_comment, '''') as `Diagnosis Comment`,
coalesce(d.disease_phase, '''') as `Disease Phase`,
coalesce(d.anatomic_site, '''') as `Anatomic Site`,
case d.age_at_diagnosis when -999 then ''Not Reported'' else coalesce(d.age_at_diagnosis, '''') end as `Age at Diagnosis (days)`,
coalesce(fo.vital_status, '''') as `Vital Status`
Order by p.participant_id limit 100'

$ws.Range('B6').Value2 = 'CALL {
MATCH (file:clinical_measure_file)
        MATCH (p:participant)-[:of_clinical_measure_file]-(file)
        MATCH (st:study)<-[:of_participant]-(p)
        Where st.phs_accession in [''phs003111'']
        with file, p, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, p, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, p, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Clinical data'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        p.participant_id as participant_id,
        "" AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file:clinical_measure_file)
        MATCH (st:study)<-[:of_clinical_measure_file]-(file)
        Where st.phs_accession in [''phs003111'']
        with file, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Clinical data'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        "" as participant_id,
        "" AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file:methylation_array_file)
        MATCH (p:participant)<-[:of_sample]-(sm1:sample)<-[*0..2]-(sm:sample)<-[:of_methylation_array_file]-(file)
        MATCH (st:study)<-[:of_participant]-(p)
        Where st.phs_accession in [''phs003111'']
        with file, p, sm1, sm, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, p, sm1, sm, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, p, sm1, sm, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Methylation array'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        p.participant_id as participant_id,
        CASE sm1.sample_id WHEN sm.sample_id THEN sm.sample_id
                    ELSE sm1.sample_id + '','' + sm.sample_id END AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file:pathology_file)
        MATCH (p:participant)<-[:of_sample]-(sm1:sample)<-[*0..2]-(sm:sample)<-[:of_pathology_file]-(file)
        MATCH (st:study)<-[:of_participant]-(p)
        Where st.phs_accession in [''phs003111'']
        with file, p, sm1, sm, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, p, sm1, sm, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, p, sm1, sm, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Pathology imaging'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        p.participant_id as participant_id,
        CASE sm1.sample_id WHEN sm.sample_id THEN sm.sample_id
                    ELSE sm1.sample_id + '','' + sm.sample_id END AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file:radiology_file)
        MATCH (p:participant)<-[:of_radiology_file]-(file)
        MATCH (st:study)<-[:of_participant]-(p)
        Where st.phs_accession in [''phs003111'']
        with file, p, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, p, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, p, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Radiology imaging'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        p.participant_id as participant_id,
        "" AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file:single_cell_sequencing_file)
        MATCH (p:participant)<-[:of_sample]-(sm1:sample)<-[*0..2]-(sm:sample)<-[:of_single_cell_sequencing_file]-(file)
        MATCH (st:study)<-[:of_participant]-(p)
        Where st.phs_accession in [''phs003111'']
        with file, p, sm1, sm, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, p, sm1, sm, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, p, sm1, sm, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Single Cell Sequencing'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        p.participant_id as participant_id,
        CASE sm1.sample_id WHEN sm.sample_id THEN sm.sample_id
                    ELSE sm1.sample_id + '','' + sm.sample_id END AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file:sequencing_file)
        MATCH (p:participant)<-[:of_sample]-(sm1:sample)<-[*0..2]-(sm:sample)<-[:of_sequencing_file]-(file)
        MATCH (st:study)<-[:of_participant]-(p)
        Where st.phs_accession in [''phs003111'']
        with file, p, sm1, sm, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, p, sm1, sm, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, p, sm1, sm, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Sequencing'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        p.participant_id as participant_id,
        CASE sm1.sample_id WHEN sm.sample_id THEN sm.sample_id
                    ELSE sm1.sample_id + '','' + sm.sample_id END AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file:cytogenomic_file)
        MATCH (p:participant)<-[:of_sample]-(sm1:sample)<-[*0..2]-(sm:sample)<-[:of_cytogenomic_file]-(file)
        MATCH (st:study)<-[:of_participant]-(p)
        Where st.phs_accession in [''phs003111'']
        with file, p, sm1, sm, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, p, sm1, sm, st,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, p, sm1, sm, st, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        ''Cytogenomic'' AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        p.participant_id as participant_id,
        CASE sm1.sample_id WHEN sm.sample_id THEN sm.sample_id
                    ELSE sm1.sample_id + '','' + sm.sample_id END AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
        UNION ALL
        MATCH (file)
        WHERE (file:sequencing_file OR file:pathology_file OR file:methylation_array_file OR file:single_cell_sequencing_file OR file:cytogenomic_file)
        MATCH (st:study)<-[:of_cell_line|of_pdx]-(cl)<--(sm:sample)
        Where (cl: cell_line or cl: pdx) and st.phs_accession in [''phs003111'']
        MATCH (sm)<--(file)
        with file, sm, st,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(file.file_size)/log(1024))) as i,
        2 as precision
        WITH file, st, sm,
        file.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
        WITH file, st, sm, unit,
        round(factor * value)/factor AS size
        RETURN DISTINCT
        file.file_name AS file_name,
        CASE LABELS(file)[0]
                WHEN ''sequencing_file'' THEN ''Sequencing''
                WHEN ''single_cell_sequencing_file'' THEN ''Single Cell Sequencing''
                WHEN ''cytogenomic_file'' THEN ''Cytogenomic''
                WHEN ''pathology_file'' THEN ''Pathology imaging''
                WHEN ''methylation_array_file'' THEN ''Methylation array'' END AS file_category,
        file.file_description AS file_description,
        file.file_type AS file_type,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS file_size,
        st.study_id AS study_id,
        "" as participant_id,
        sm.sample_id AS sample_id,
        file.dcf_indexd_guid AS guid,
        file.md5sum AS md5sum
}
RETURN file_name AS `File Name`,
file_category As `File Category`,
file_description As `File Description`,
file_type As `File Type`,
file_size As `File Size`,
study_id As `Study ID`,
participant_id As `Participant ID`,
sample_id As `Sample ID`,
guid As `GUID`,
md5sum As `MD5sum`
ORDER BY file_name LIMIT 100'

# Move the saved selection/active cell on the startup sheet.
$ws.Range('C2').Select()
